$d = $word.ActiveDocument

$replacements = @(
    @("50×99=", "45×36="),
    @("86×59=", "79×25="),
    @("46×11=", "57×60="),
    @("23×70=", "99×69="),
    @("36×49=", "79×16="),
    @("23×89=", "47×36="),
    @("37×49=", "43×93="),
    @("28×20=", "14×13="),
    @("88×16=", "63×54="),
    @("40×42=", "14×82="),
    @("36×84=", "61×90="),
    @("40×25=", "85×67="),
    @("36×93=", "80×30="),
    @("40×58=", "41×97="),
    @("98×87=", "79×66="),
    @("76×62=", "17×97="),
    @("32×72=", "64×52="),
    @("95×58=", "59×83="),
    @("50×83=", "84×26="),
    @("87×96=", "20×53="),
    @("85×54=", "62×26="),
    @("83×85=", "94×65="),
    @("79×53=", "80×51="),
    @("90×73=", "75×13="),
    @("68×81=", "99×30=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
